$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "02-12-2022"
$ws.Range("L2").Value = "58571703"
